# Generate Report for Handoff
# This script updates the localization-status workbook so that the two
# source documents (previously "d60ccee0..." and "da9d4636...") are
# replaced by a new pair of documents ("6d9ff285..." and "ffff4aff7191...")
# that are freshly queued for handoff (not yet handed back).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New identifiers used throughout the report
# ---------------------------------------------------------------------
$docA      = "6d9ff285-2ced-47f6-92c3-438f7494b03e"
$docB      = "ffff4aff7191-ec5f-4eb0-9031-30673ea919ee"
$xlfHash   = "ca8f6f38853e99d0911bf5b6652330329e7ac047"

$docA_md       = "$docA.md"
$docB_md       = "$docB.md"
$docA_zh_xlf   = "$docA.$xlfHash.zh-cn.xlf"
$docA_de_xlf   = "$docA.$xlfHash.de-de.xlf"

$status        = "Ready for handoff"
$overviewDate  = "2016-52-19 08:52:04"
$zhHandoffDate = "2016-03-19 08:52:00"
$deHandoffDate = "2016-03-19 08:52:04"
$neverDate     = "0001-01-01 00:00:00"

# Hyperlink target URL bases (kept consistent with the existing pattern,
# only the file identifiers change)
$mdBase      = "https://github.com/OpenLocalizationTest/oltest/blob/cb1a63928fc863f9190ebdce8053f33cdce2a0ae/e2e"
$zhOffBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4edfde31e72067446eed6d89ce57dfb0722b4171/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deOffBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9480dac98c276b8ec77f13fd79e995312b6364f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = $docA_md
$ov.Range("B2").Value = $status
$ov.Range("C2").Value = $status
$ov.Range("D2").Value = $overviewDate

$ov.Range("A3").Value = $docB_md
$ov.Range("B3").Value = $status
$ov.Range("C3").Value = $status
$ov.Range("D3").Value = $overviewDate

$ov.Hyperlinks.Add($ov.Range("A2"), "$mdBase/$docA_md", "", "", $docA_md) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "$mdBase/$docB_md", "", "", $docB_md) | Out-Null

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = $docA_md
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $status
$zh.Range("D2").Value = $docA_zh_xlf
$zh.Range("E2").Value = $zhHandoffDate
$zh.Range("F2:G2").Clear()
$zh.Range("H2").Value = $neverDate
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = $docB_md
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $status
$zh.Range("D3").Value = $docA_zh_xlf
$zh.Range("E3").Value = $zhHandoffDate
$zh.Range("F3:G3").Clear()
$zh.Range("H3").Value = $neverDate
$zh.Range("I3").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A2"), "$mdBase/$docA_md", "", "", $docA_md) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), "$mdBase/$docA_md", "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "$zhOffBase/$docA_zh_xlf", "", "", $docA_zh_xlf) | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), "$mdBase/$docB_md", "", "", $docB_md) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), "$mdBase/$docB_md", "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "$zhOffBase/$docA_zh_xlf", "", "", $docA_zh_xlf) | Out-Null

# =======================================================================
# Sheet "de-de"
# =======================================================================
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

$de.Range("A2").Value = $docA_md
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $status
$de.Range("D2").Value = $docA_de_xlf
$de.Range("E2").Value = $deHandoffDate
$de.Range("F2:G2").Clear()
$de.Range("H2").Value = $neverDate
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = $docB_md
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $status
$de.Range("D3").Value = $docA_de_xlf
$de.Range("E3").Value = $deHandoffDate
$de.Range("F3:G3").Clear()
$de.Range("H3").Value = $neverDate
$de.Range("I3").Value = "Include"

$de.Hyperlinks.Add($de.Range("A2"), "$mdBase/$docA_md", "", "", $docA_md) | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), "$mdBase/$docA_md", "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "$deOffBase/$docA_de_xlf", "", "", $docA_de_xlf) | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), "$mdBase/$docB_md", "", "", $docB_md) | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), "$mdBase/$docB_md", "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "$deOffBase/$docA_de_xlf", "", "", $docA_de_xlf) | Out-Null
